# Update handback status timestamps (regenerate report for handback)
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-16 23:04:17"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-16 23:04:12"
$wsZhCn.Range("K2").Value = "2016-08-16 23:04:29"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-16 23:04:36"
